$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) figures for rows 3,5,7,8,10,11,12
# on both the "展览" and "全部类型" sheets (they mirror the same data).
$updates = @{
    3  = 6289
    5  = 21
    7  = 1897
    8  = 1431
    10 = 957
    11 = 253
    12 = 5586
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
